$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain plain text (avoids Excel
# auto-converting numeric-looking strings like "25.30" or "1.001" into numbers,
# which would drop significant trailing/leading zeros). We temporarily switch the
# cell to Text format, assign the value, then restore the cell's original style so no
# visual/style residue is left behind.
function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = $origStyle
}

Set-TextValue "D2" '29.843.71'
$ws.Range("E2").Value = '  -0.90%  '
Set-TextValue "D3" '1.893.60'
$ws.Range("E3").Value = '  -0.57%  '
Set-TextValue "D4" '1.001'
$ws.Range("E4").Value = '  +0.16%  '
Set-TextValue "D5" '0.7907'
$ws.Range("E5").Value = '  -5.68%  '
Set-TextValue "D6" '243.34'
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E8").Value = '  -4.43%  '
Set-TextValue "D9" '25.30'
$ws.Range("E9").Value = '  -5.29%  '
Set-TextValue "D10" '0.07222'
$ws.Range("E10").Value = '  +2.03%  '
Set-TextValue "D11" '0.08094'
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D12" '0.7648'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D13" '5.545'
$ws.Range("E13").Value = '  +5.01%  '
Set-TextValue "D14" '1.913.75'
$ws.Range("E14").Value = '  +0.51%  '
Set-TextValue "D15" '92.42'
$ws.Range("E15").Value = '  -0.30%  '
Set-TextValue "D16" '6.149'
$ws.Range("E16").Value = '  +4.45%  '
Set-TextValue "D17" '29.848.50'
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("E18").Value = '  -2.01%  '
Set-TextValue "D19" '243.74'
$ws.Range("E19").Value = '  -0.50%  '
Set-TextValue "D20" '0.000007770'
$ws.Range("E20").Value = '  -0.09%  '
Set-TextValue "D21" '2.157.75'
$ws.Range("E21").Value = '  +0.08%  '
Set-TextValue "D22" '1.001'
$ws.Range("E22").Value = '  +0.10%  '
Set-TextValue "D23" '8.122'
$ws.Range("E23").Value = '  +15.76%  '
Set-TextValue "D24" '1.001'
$ws.Range("E24").Value = '  +0.17%  '
Set-TextValue "D25" '0.1649'
$ws.Range("E25").Value = '  -7.14%  '
Set-TextValue "D26" '9.376'
$ws.Range("E26").Value = '  +0.74%  '
Set-TextValue "D27" '163.50'
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("E28").Value = '  -1.57%  '
Set-TextValue "D29" '2.048'
$ws.Range("E29").Value = '  -2.58%  '
$ws.Range("E30").Value = '  +2.50%  '
Set-TextValue "D31" '1.547'
$ws.Range("E31").Value = '  +1.72%  '
Set-TextValue "D32" '4.464'
$ws.Range("E32").Value = '  +3.78%  '
Set-TextValue "D33" '4.093'
$ws.Range("E33").Value = '  +0.16%  '
Set-TextValue "D34" '0.05539'
$ws.Range("E34").Value = '  -6.92%  '
$ws.Range("E35").Value = '  -0.43%  '
Set-TextValue "D36" '0.7381'
$ws.Range("E36").Value = '  +0.54%  '
Set-TextValue "D37" '0.9999'
$ws.Range("E37").Value = '  +0.17%  '
Set-TextValue "D38" '2.623'
$ws.Range("E38").Value = '  -3.18%  '
Set-TextValue "D39" '0.01925'
$ws.Range("E39").Value = '  +0.06%  '
Set-TextValue "D40" '2.779'
$ws.Range("E40").Value = '  -0.29%  '
Set-TextValue "D41" '1.147.82'
$ws.Range("E41").Value = '  +14.41%  '
Set-TextValue "D42" '73.85'
$ws.Range("E42").Value = '  +1.06%  '
Set-TextValue "D43" '0.4421'
Set-TextValue "D44" '5.872'
$ws.Range("E44").Value = '  -1.47%  '
Set-TextValue "D45" '0.8493'
$ws.Range("E45").Value = '  -1.01%  '
Set-TextValue "D46" '104.21'
$ws.Range("E46").Value = '  +2.05%  '
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D48" '1.874'
$ws.Range("E48").Value = '  -1.84%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D49" '9.992'
$ws.Range("E49").Value = '  +1.59%  '
Set-TextValue "D50" '7.439'
$ws.Range("E50").Value = '  -1.90%  '
Set-TextValue "D51" '2.996'
$ws.Range("E51").Value = '  +9.85%  '
